# Adds in row label to dict
#
# This reproduces the source diff:
#   - Sheet1!A2 gets the new label "row number"
#   - Sheet2!A5 gets the same new label "row number"
#     (both pull the same new shared string, appended once to sharedStrings.xml)
#   - Sheet1's sheetView selection becomes A1:I5 (whole header+data block)
#   - Sheet2's sheetView selection/active cell moves to G22

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Add the new "row number" label cells -------------------------------
$ws1.Range("A2").Value = "row number"
$ws2.Range("A5").Value = "row number"

# --- Update the on-screen selections for each sheet ----------------------
$ws1.Activate()
$ws1.Range("A1:I5").Select()

$ws2.Activate()
$ws2.Range("G22").Select()
